# Fix the DOCUSIGN_INC workbook: the open/close/high/low price, shares
# outstanding and fixed_ticker columns (D:I) on Sheet1 had been polluted
# with data copied from other tickers' files. Restore them to DOCU's own
# values, which also drops the now-unused extra ticker strings (SNPS, IBM,
# PRO, G, INTC, ZS, FICO, BABA, SMCI, GOOGL, MRVL, UBER, MU, NFLX, BIDU,
# LPSN, ORCL, BAH, AMZN, CDNS, CAN, DDOG, CGNT, LSCC, HUBS, TSLA, ACN, APP,
# AI, PTC) from the shared-string table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row -> open, close, high, low, shares_outstanding
$data = @{
    2  = @(38,                  38.63000106811523,  40.88999938964844,  37,                  201104117)
    3  = @(38,                  38.63000106811523,  40.88999938964844,  37,                  201104117)
    4  = @(38,                  38.63000106811523,  40.88999938964844,  37,                  201104117)
    5  = @(38,                  38.63000106811523,  40.88999938964844,  37,                  201104117)
    6  = @(38,                  38.63000106811523,  40.88999938964844,  37,                  201104117)
    7  = @(38.70000076293945,   49.81000137329102,  51.20100021362305,  37.84999847412109,   201104117)
    8  = @(54.29999923706055,   62.43999862670898,  68.34999847412109,  54.19499969482422,   201104117)
    9  = @(42.13999938964844,   41.7599983215332,    45.15000152587891,  35.06000137329102,   201104117)
    10 = @(49.33000183105469,   55.15000152587891,  56.22999954223633,  48.20999908447266,   201104117)
    11 = @(57.09999847412109,   56.06000137329102,  57.41999816894531,  50.02000045776367,   201104117)
    12 = @(51.72999954223633,   46.68999862670898,  53,                 43.13000106811523,   201104117)
    13 = @(66.77999877929688,   71.20999908447266,  73.22000122070312,  64.25,               201104117)
    14 = @(79.09999847412109,   86.30999755859375,  92.5500030517578,   76.18000030517578,   201104117)
    15 = @(103.1900024414062,   139.7400054931641,  141.4299926757812,  101.879997253418,    201104117)
    16 = @(218.7200012207031,   223,                229.8300018310547,  189.1199951171875,   201104117)
    17 = @(204.2100067138672,   227.8800048828125,  246,                185.3589935302734,   201104117)
    18 = @(233.25,              226.6600036621093,  275.4599914550781,  219.1499938964844,   201104117)
    19 = @(224.0299987792969,   201.6199951171875,  224.0299987792969,  179.4900054931641,   201104117)
    20 = @(299.7999877929688,   296.239990234375,    314.760009765625,   281.2799987792969,   201104117)
    21 = @(277.2699890136719,   246.3600006103516,  288.135009765625,   233.3300018310547,   201104117)
    22 = @(127.5899963378906,   118.4300003051758,  131.9069976806641,  100,                 201104117)
    23 = @(80.73000335693359,   83.91000366210938,  88.58000183105469,  64.83999633789062,   201104117)
    24 = @(63.31000137329102,   58.22000122070312,  77.23999786376953,  56.91999816894531,   201104117)
    25 = @(49.72999954223633,   47.06999969482422,  54.32500076293945,  39.56999969482422,   201104117)
    26 = @(61,                  61.34999847412109,  69.44999694824219,  58.06000137329102,   201104117)
    27 = @(49.20000076293945,   56.40000152587891,  56.97999954223633,  46.79999923706055,   201104117)
    28 = @(53.36999893188477,   50.29999923706055,  53.52000045776367,  46.65999984741211,   201104117)
    29 = @(38.81000137329102,   43.09999847412109,  44.33000183105469,  38.20000076293945,   201104117)
    30 = @(61.15000152587891,   53.27000045776367,  62.68999862670898,  49.11999893188477,   201104117)
    31 = @(56.36000061035156,   54.7400016784668,    61.13999938964844,  53.5,                201104117)
    32 = @(55.43999862670898,   59.20999908447266,  59.77999877929688,  48.79999923706055,   201104117)
    33 = @(69.87000274658203,   79.69000244140625,  86.75,              69.27999877929688,   201104117)
    34 = @(93.06199645996094,   83.16999816894531,  99.3000030517578,   78.51399993896484,   201104117)
    35 = @(82.40000152587891,   88.61000061035156,  92.97000122070312,  80.54000091552734,   201104117)
    36 = @(74.88999938964844,   76.66000366210938,  77.19000244140625,  66.35199737548828,   201104117)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 4).Value = $vals[0]   # D open_price
    $ws.Cells.Item($row, 5).Value = $vals[1]   # E close_price
    $ws.Cells.Item($row, 6).Value = $vals[2]   # F high_price
    $ws.Cells.Item($row, 7).Value = $vals[3]   # G low_price
    $ws.Cells.Item($row, 8).Value = $vals[4]   # H shares_outstanding
    $ws.Cells.Item($row, 9).Value = "DOCU"     # I fixed_ticker
}
